# Updates the cryptocurrency price (column D) and Volume(1h) (column E)
# values on the active worksheet to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.723.37"
$ws.Range("E2").Value = "  +1.87%  "

$ws.Range("D3").Value = "1.871.04"
$ws.Range("E3").Value = "  +1.87%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9987"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.75"
$ws.Range("E5").Value = "  +1.54%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9991"
$ws.Range("E6").Value = "  -0.13%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4788"
$ws.Range("E7").Value = "  +2.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2838"
$ws.Range("E8").Value = "  +5.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06528"
$ws.Range("E9").Value = "  +4.17%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.71"
$ws.Range("E10").Value = "  +16.79%  "

$ws.Range("D11").Value = "1.912.84"
$ws.Range("E11").Value = "  +4.06%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07498"
$ws.Range("E12").Value = "  +1.28%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "94.60"
$ws.Range("E13").Value = "  +13.03%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.099"
$ws.Range("E14").Value = "  +3.68%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6523"
$ws.Range("E15").Value = "  +5.48%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "297.54"
$ws.Range("E16").Value = "  +30.81%  "

$ws.Range("D17").Value = "30.699.38"
$ws.Range("E17").Value = "  +2.05%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9984"
$ws.Range("E18").Value = "  -0.18%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.01"
$ws.Range("E19").Value = "  +5.42%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007480"
$ws.Range("E20").Value = "  +2.82%  "

$ws.Range("D21").Value = "2.117.88"
$ws.Range("E21").Value = "  +1.73%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9997"
$ws.Range("E22").Value = "  -0.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.185"
$ws.Range("E23").Value = "  +6.65%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.108"
$ws.Range("E24").Value = "  +4.79%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "168.94"
$ws.Range("E25").Value = "  +2.29%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.247"
$ws.Range("E26").Value = "  +0.70%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.51"
$ws.Range("E27").Value = "  +10.23%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.956"
$ws.Range("E28").Value = "  +4.49%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1046"
$ws.Range("E29").Value = "  +1.76%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.348"
$ws.Range("E30").Value = "  -1.60%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.120"
$ws.Range("E31").Value = "  +1.00%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.962"
$ws.Range("E32").Value = "  +4.70%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04965"
$ws.Range("E33").Value = "  +3.49%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.177"
$ws.Range("E34").Value = "  +3.71%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7197"
$ws.Range("E35").Value = "  +1.75%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.707"
$ws.Range("E36").Value = "  +0.25%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01939"
$ws.Range("E37").Value = "  +4.24%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.708"
$ws.Range("E38").Value = "  +2.12%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.037"
$ws.Range("E39").Value = "  +5.64%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8880"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "107.12"
$ws.Range("E41").Value = "  +2.74%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9985"
$ws.Range("E42").Value = "  -0.27%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4197"
$ws.Range("E43").Value = "  +5.00%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.570"
$ws.Range("E44").Value = "  +0.78%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.349"
$ws.Range("E45").Value = "  +5.83%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.82"
$ws.Range("E46").Value = "  +8.63%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1228"
$ws.Range("E47").Value = "  +3.23%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "34.71"
$ws.Range("E48").Value = "  +6.49%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.745"
$ws.Range("E49").Value = "  +2.35%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.393"
$ws.Range("E50").Value = "  +2.72%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05543"
$ws.Range("E51").Value = "  +0.66%  "

